$wb = $excel.ActiveWorkbook

# --- ALC row 12 (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 280
$ws.Range("I12").Value = 293.33334
$ws.Range("J12").Value = 240
$ws.Range("K12").Value = 293.33334
$ws.Range("L12").Value = 240
$ws.Range("M12").Value = -123.33334
$ws.Range("N12").Value = -580

# --- ALC row 132 (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2370.625
$ws.Range("I132").Value = 2461.842
$ws.Range("J132").Value = 2024
$ws.Range("K132").Value = 7385.526
$ws.Range("L132").Value = 6072
$ws.Range("M132").Value = -4855.526
$ws.Range("N132").Value = -11132

# --- ALC row 137 (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1034.1786
$ws.Range("I137").Value = 775.36365
$ws.Range("K137").Value = 2326.09095
$ws.Range("M137").Value = 223.9090500000002

# --- ALC row 140 (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 71431.11
$ws.Range("J140").Value = 71431.11
$ws.Range("L140").Value = 71431.11
$ws.Range("N140").Value = -81791.11

# --- ALC row 141 (hunk 4) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1569.4736
$ws.Range("I141").Value = 1301.25
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 3903.75
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = 1276.25
$ws.Range("N141").Value = -19360

# --- ARM row 3 (hunk 5) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -885

# --- ARM row 61 (hunk 6) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2009.1389
$ws.Range("I61").Value = 1928.6072
$ws.Range("J61").Value = 2291
$ws.Range("K61").Value = 1928.6072
$ws.Range("L61").Value = 2291
$ws.Range("M61").Value = -1716.6072
$ws.Range("N61").Value = -2715

# --- ARM row 74 (hunk 7) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3039.5334
$ws.Range("I74").Value = 3240.7273
$ws.Range("K74").Value = 3240.7273
$ws.Range("M74").Value = -2366.7273

# --- ARM row 77 (hunk 8) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3039.5334
$ws.Range("I77").Value = 3240.7273
$ws.Range("K77").Value = 16203.6365
$ws.Range("M77").Value = -11835.6365

# --- ARM row 132 (hunk 9) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2545.64
$ws.Range("I132").Value = 2106.0527
$ws.Range("K132").Value = 6318.158100000001
$ws.Range("M132").Value = -3788.158100000001

# --- ARM row 136 (hunk 10) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2009.1389
$ws.Range("I136").Value = 1928.6072
$ws.Range("J136").Value = 2291
$ws.Range("K136").Value = 5785.821599999999
$ws.Range("L136").Value = 6873
$ws.Range("M136").Value = -3235.821599999999
$ws.Range("N136").Value = -11973

# --- BSM row 7 (hunk 11) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 850
$ws.Range("I7").Value = 750
$ws.Range("J7").Value = 950
$ws.Range("K7").Value = 750
$ws.Range("L7").Value = 950
$ws.Range("M7").Value = -637
$ws.Range("N7").Value = -1176

# --- BSM row 107 (hunk 12) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1430478.6
$ws.Range("I107").Value = 1793.6666
$ws.Range("J107").Value = 2501992.2
$ws.Range("K107").Value = 1793.6666
$ws.Range("L107").Value = 2501992.2
$ws.Range("M107").Value = 126.3334
$ws.Range("N107").Value = -2505832.2

# --- BSM row 134 (hunk 13) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4234.3335
$ws.Range("I134").Value = 981.5854
$ws.Range("J134").Value = 37575
$ws.Range("K134").Value = 2944.7562
$ws.Range("L134").Value = 112725
$ws.Range("M134").Value = -409.7562000000003
$ws.Range("N134").Value = -117795

# --- BSM row 140 (hunk 14) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 86840
$ws.Range("J140").Value = 86840
$ws.Range("L140").Value = 86840
$ws.Range("N140").Value = -97200

# --- CRP row 2 (hunk 15) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 240
$ws.Range("I2").Value = 240
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 240
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -127
$ws.Range("N2").ClearContents()

# --- CRP row 31 (hunk 16) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2572.6428
$ws.Range("J31").Value = 5241.4
$ws.Range("L31").Value = 5241.4
$ws.Range("N31").Value = -5831.4

# --- CRP row 34 (hunk 17) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2572.6428
$ws.Range("J34").Value = 5241.4
$ws.Range("L34").Value = 5241.4
$ws.Range("N34").Value = -5645.4

# --- CRP row 58 (hunk 18) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1552
$ws.Range("I58").Value = 1259.4286
$ws.Range("J58").Value = 1893.3334
$ws.Range("K58").Value = 1259.4286
$ws.Range("L58").Value = 1893.3334
$ws.Range("M58").Value = -1056.4286
$ws.Range("N58").Value = -2299.3334

# --- CRP row 136 (hunk 19) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1552
$ws.Range("I136").Value = 1259.4286
$ws.Range("J136").Value = 1893.3334
$ws.Range("K136").Value = 3778.2858
$ws.Range("L136").Value = 5680.0002
$ws.Range("M136").Value = -1228.2858
$ws.Range("N136").Value = -10780.0002

# --- CUL row 131 (hunk 20) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 812.21
$ws.Range("I131").Value = 625
$ws.Range("J131").Value = 833.0111000000001
$ws.Range("K131").Value = 1875
$ws.Range("L131").Value = 2499.0333
$ws.Range("M131").Value = 3165
$ws.Range("N131").Value = -12579.0333

# --- GSM row 5 (hunk 21) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 38750
$ws.Range("I5").Value = 5000
$ws.Range("K5").Value = 5000
$ws.Range("M5").Value = -4888

# --- GSM row 21 (hunk 22) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 3334166.8
$ws.Range("J21").Value = 1250
$ws.Range("L21").Value = 1250
$ws.Range("N21").Value = -1596

# --- GSM row 30 (hunk 23) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 3334166.8
$ws.Range("J30").Value = 1250
$ws.Range("L30").Value = 1250
$ws.Range("N30").Value = -1460

# --- GSM row 140 (hunk 24) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 27652
$ws.Range("J140").Value = 27652
$ws.Range("L140").Value = 27652
$ws.Range("N140").Value = -38012

# --- LTW row 2 (hunk 25) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6000000
$ws.Range("J2").Value = 6000000
$ws.Range("L2").Value = 6000000
$ws.Range("N2").Value = -6000224

# --- LTW row 40 (hunk 26) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2398.7896
$ws.Range("I40").Value = 1980
$ws.Range("K40").Value = 1980
$ws.Range("M40").Value = -1844

# --- LTW row 122 (hunk 27) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2786.0908
$ws.Range("I122").Value = 2213.25
$ws.Range("J122").Value = 4313.6665
$ws.Range("K122").Value = 6639.75
$ws.Range("L122").Value = 12940.9995
$ws.Range("M122").Value = -4189.75
$ws.Range("N122").Value = -17840.9995

# --- LTW row 132 (hunk 28) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5268.0527
$ws.Range("I132").Value = 6099.923
$ws.Range("J132").Value = 3465.6667
$ws.Range("K132").Value = 18299.769
$ws.Range("L132").Value = 10397.0001
$ws.Range("M132").Value = -15769.769
$ws.Range("N132").Value = -15457.0001

# --- LTW row 139 (hunk 29) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 79550
$ws.Range("J139").Value = 79550
$ws.Range("L139").Value = 79550
$ws.Range("N139").Value = -89830

# --- WVR row 2 (hunk 30) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2466.6667
$ws.Range("J2").Value = 2800
$ws.Range("L2").Value = 2800
$ws.Range("N2").Value = -3024

# --- WVR row 4 (hunk 31) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 600.6667
$ws.Range("I4").Value = 401.33334
$ws.Range("J4").Value = 800
$ws.Range("K4").Value = 401.33334
$ws.Range("L4").Value = 800
$ws.Range("M4").Value = -288.33334
$ws.Range("N4").Value = -1026

# --- WVR row 15 (hunk 32) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

# --- WVR row 19 (hunk 33) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 1724.8
$ws.Range("J19").Value = 1729.75
$ws.Range("L19").Value = 1729.75
$ws.Range("N19").Value = -2077.75

# --- WVR row 138 (hunk 34) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 48533.332
$ws.Range("J138").Value = 48533.332
$ws.Range("L138").Value = 48533.332
$ws.Range("N138").Value = -58813.332

# --- WVR row 139 (hunk 35) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 38575
$ws.Range("J139").Value = 38575
$ws.Range("L139").Value = 38575
